# This script applies corrected stock-quantity/value figures to the
# CryCompanywiseStockReport worksheet, matching a re-run of the report
# (quantities/values recalculated for a number of items, two pairs of
# duplicate-named rows re-ordered, and the affected Sub Total / Grand
# Total rows updated accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F21").Value = 150
$ws.Range("G21").Value = 3850.5
$ws.Range("F28").Value = 70
$ws.Range("G28").Value = 3122
$ws.Range("B32").Value = 12939.22
$ws.Range("F135").Value = 27
$ws.Range("G135").Value = 837.8099999999999
$ws.Range("B138").Value = 2668.95
$ws.Range("F144").Value = 1144
$ws.Range("G144").Value = 9666.799999999999
$ws.Range("F145").Value = 568
$ws.Range("G145").Value = 4538.32
$ws.Range("B147").Value = 16983.39
$ws.Range("F150").Value = 41
$ws.Range("G150").Value = 1906.09
$ws.Range("B156").Value = 33901.54
$ws.Range("B227").Value = 63520
$ws.Range("E227").Value = 153.4
$ws.Range("F227").Value = 67
$ws.Range("G227").Value = 9666.76
$ws.Range("B228").Value = 55373
$ws.Range("E228").Value = 163.62
$ws.Range("F228").Value = -94
$ws.Range("G228").Value = -13562.32
$ws.Range("B229").Value = 63531
$ws.Range("E229").Value = 152.53
$ws.Range("F229").Value = 65
$ws.Range("G229").Value = 9326.200000000001
$ws.Range("B230").Value = 57802
$ws.Range("E230").Value = 162.71
$ws.Range("F230").Value = -79
$ws.Range("G230").Value = -11334.92
$ws.Range("F282").Value = 5
$ws.Range("G282").Value = 268.5
$ws.Range("F294").Value = 39
$ws.Range("G294").Value = 2783.04
$ws.Range("B304").Value = 186385.07
$ws.Range("B322").Value = 47097
$ws.Range("D322").Value = 112.28
$ws.Range("E322").Value = 134.16
$ws.Range("F322").Value = 15
$ws.Range("G322").Value = 1684.2
$ws.Range("B323").Value = 58047
$ws.Range("D323").Value = 105.54
$ws.Range("E323").Value = 126.1
$ws.Range("F323").Value = 40
$ws.Range("G323").Value = 4221.6
$ws.Range("F338").Value = 80
$ws.Range("G338").Value = 1896
$ws.Range("F345").Value = 66
$ws.Range("G345").Value = 4053.06
$ws.Range("B346").Value = 27391.86
$ws.Range("B364").Value = 53602
$ws.Range("E364").Value = 15.69
$ws.Range("F364").Value = -231
$ws.Range("G364").Value = -3037.65
$ws.Range("B365").Value = 65068
$ws.Range("E365").Value = 13.97
$ws.Range("F365").Value = 63
$ws.Range("G365").Value = 828.45
$ws.Range("B382").Value = 45702
$ws.Range("E382").Value = 31.43
$ws.Range("F382").Value = -215
$ws.Range("G382").Value = -5654.5
$ws.Range("B383").Value = 64919
$ws.Range("E383").Value = 27.97
$ws.Range("F383").Value = 61
$ws.Range("G383").Value = 1604.3
$ws.Range("B385").Value = 65067
$ws.Range("E385").Value = 15.65
$ws.Range("F385").Value = 126
$ws.Range("G385").Value = 1855.98
$ws.Range("B386").Value = 53595
$ws.Range("E386").Value = 17.61
$ws.Range("F386").Value = -335
$ws.Range("G386").Value = -4934.55
$ws.Range("F455").Value = 48
$ws.Range("G455").Value = 3053.28
$ws.Range("B460").Value = 14243.06
$ws.Range("B463").Value = 60025
$ws.Range("E463").Value = 37.22
$ws.Range("F463").Value = -98
$ws.Range("G463").Value = -3217.34
$ws.Range("B464").Value = 64833
$ws.Range("E464").Value = 34.9
$ws.Range("F464").Value = 95
$ws.Range("G464").Value = 3118.85
$ws.Range("F477").Value = 12
$ws.Range("G477").Value = 544.08
$ws.Range("B478").Value = 544.08
$ws.Range("F485").Value = 22
$ws.Range("G485").Value = 3860.34
$ws.Range("B488").Value = 31798.34
$ws.Range("F492").Value = 61
$ws.Range("G492").Value = 7963.55
$ws.Range("B493").Value = 12983.72
$ws.Range("F550").Value = 3
$ws.Range("G550").Value = 244.68
$ws.Range("F555").Value = 32
$ws.Range("G555").Value = 2225.92
$ws.Range("F556").Value = 2
$ws.Range("G556").Value = 229.72
$ws.Range("B560").Value = 6607.71
$ws.Range("B619").Value = 1884417.83
$ws.Range("B620").Value = 1884417.83
